$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new finished match (id 212, B=6788924) is inserted at sheet row 214,
# pushing the existing rows 214-216 down to 215-217 (with several of their
# odds columns recalculated). Row 217 is brand-new (was previously absent).
#
# Strategy: first clone the A/E cell styles (bold+border / date format) for
# the brand-new row 217 from row 216 (style-only, single cell, so no stray
# cells / new style entries get created), then overwrite every data cell
# for rows 214-217 with the exact final values from the target workbook.
# ---------------------------------------------------------------------------

# Bring style forward onto the new row 217 for the two styled columns only.
$ws.Range("A216").Copy($ws.Range("A217"))
$ws.Range("E216").Copy($ws.Range("E217"))

# ---- Row 214 (new match: NK Varazdin vs Hajduk Split, 1-1 draw) ----------
$ws.Range("A214").Value2 = 212
$ws.Range("B214").Value2 = 6788924
$ws.Range("C214").Value2 = "Croatia HNL"
$ws.Range("D214").Value2 = "Croatia HNL"
$ws.Range("E214").Value2 = 45346.45833333334
$ws.Range("F214").Value2 = "NK Varazdin"
$ws.Range("G214").Value2 = "Hajduk Split"
$ws.Range("H214").Value2 = 1
$ws.Range("I214").Value2 = 1
$ws.Range("J214").Value2 = "D"
$ws.Range("K214").Value2 = 4.5
$ws.Range("L214").Value2 = 3.8
$ws.Range("M214").Value2 = 1.727
$ws.Range("N214").Value2 = 6.5
$ws.Range("O214").Value2 = 4.2
$ws.Range("P214").Value2 = 1.5
$ws.Range("Q214").Value2 = 1
$ws.Range("R214").Value2 = 2
$ws.Range("S214").Value2 = 1.85
$ws.Range("T214").Value2 = 2.25
$ws.Range("U214").Value2 = 1.825
$ws.Range("V214").Value2 = 2.025
$ws.Range("W214").Value2 = -1
$ws.Range("X214").Value2 = 3.2
$ws.Range("Y214").Value2 = -1
$ws.Range("Z214").Value2 = 1
$ws.Range("AA214").Value2 = -1
$ws.Range("AB214").Value2 = -0.5
$ws.Range("AC214").Value2 = 0.5125

# ---- Row 215 (was row 214: Istra 1961 vs NK Rudes, recalculated) --------
$ws.Range("A215").Value2 = 213
$ws.Range("B215").Value2 = 6769302
$ws.Range("C215").Value2 = "Croatia HNL"
$ws.Range("D215").Value2 = "Croatia HNL"
$ws.Range("E215").Value2 = 45346.54861111111
$ws.Range("F215").Value2 = "Istra 1961"
$ws.Range("G215").Value2 = "NK Rudes"
$ws.Range("H215").Value2 = 2
$ws.Range("I215").Value2 = 1
$ws.Range("J215").Value2 = "H"
$ws.Range("K215").Value2 = 1.444
$ws.Range("L215").Value2 = 4.4
$ws.Range("M215").Value2 = 7
$ws.Range("N215").Value2 = 1.5
$ws.Range("O215").Value2 = 4
$ws.Range("P215").Value2 = 7.5
$ws.Range("Q215").Value2 = -1
$ws.Range("R215").Value2 = 1.825
$ws.Range("S215").Value2 = 2.025
$ws.Range("T215").Value2 = 2.25
$ws.Range("U215").Value2 = 1.875
$ws.Range("V215").Value2 = 1.975
$ws.Range("W215").Value2 = 0.5
$ws.Range("X215").Value2 = -1
$ws.Range("Y215").Value2 = -1
$ws.Range("Z215").Value2 = 0
$ws.Range("AA215").Value2 = 0
$ws.Range("AB215").Value2 = 0.875
$ws.Range("AC215").Value2 = -1

# ---- Row 216 (was row 215: Dinamo Zagreb vs HNK Rijeka, recalculated) ---
$ws.Range("A216").Value2 = 214
$ws.Range("B216").Value2 = 6788923
$ws.Range("C216").Value2 = "Croatia HNL"
$ws.Range("D216").Value2 = "Croatia HNL"
$ws.Range("E216").Value2 = 45347.45833333334
$ws.Range("F216").Value2 = "Dinamo Zagreb"
$ws.Range("G216").Value2 = "HNK Rijeka"
$ws.Range("K216").Value2 = 1.8
$ws.Range("L216").Value2 = 3.75
$ws.Range("M216").Value2 = 4.2
$ws.Range("N216").Value2 = 2
$ws.Range("O216").Value2 = 3.6
$ws.Range("P216").Value2 = 3.75
$ws.Range("Q216").Value2 = -0.5
$ws.Range("R216").Value2 = 2.05
$ws.Range("S216").Value2 = 1.8
$ws.Range("T216").Value2 = 2.5
$ws.Range("U216").Value2 = 2
$ws.Range("V216").Value2 = 1.85
$ws.Range("W216").Value2 = 0
$ws.Range("X216").Value2 = 0
$ws.Range("Y216").Value2 = 0
$ws.Range("Z216").Value2 = 0
$ws.Range("AA216").Value2 = 0

# ---- Row 217 (was row 216: HNK Gorica vs Slaven Belupo, recalculated) ---
$ws.Range("A217").Value2 = 215
$ws.Range("B217").Value2 = 6788925
$ws.Range("C217").Value2 = "Croatia HNL"
$ws.Range("D217").Value2 = "Croatia HNL"
$ws.Range("E217").Value2 = 45347.5625
$ws.Range("F217").Value2 = "HNK Gorica"
$ws.Range("G217").Value2 = "Slaven Belupo"
$ws.Range("K217").Value2 = 2
$ws.Range("L217").Value2 = 3.5
$ws.Range("M217").Value2 = 3.6
$ws.Range("N217").Value2 = 2.375
$ws.Range("O217").Value2 = 3.4
$ws.Range("P217").Value2 = 2.875
$ws.Range("Q217").Value2 = -0.25
$ws.Range("R217").Value2 = 2.1
$ws.Range("S217").Value2 = 1.775
$ws.Range("T217").Value2 = 2.25
$ws.Range("U217").Value2 = 1.825
$ws.Range("V217").Value2 = 2.025
$ws.Range("W217").Value2 = 0
$ws.Range("X217").Value2 = 0
$ws.Range("Y217").Value2 = 0
$ws.Range("Z217").Value2 = 0
$ws.Range("AA217").Value2 = 0
